$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The taxon-specific data in row 3 and row 4 was swapped (same two
# observations, just listed in the opposite order). Columns A, B, E, F,
# G, H, Q, R hold the per-taxon data that differs between the rows;
# everything else (validation status, locality, dates, observers, ...)
# is identical between the two rows already, so only those columns need
# to be exchanged.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture all row-3 and row-4 values first so writing new values can't
# clobber a value that still needs to be read.
$row3Values = @{}
$row4Values = @{}
foreach ($col in $cols) {
    $row3Values[$col] = $ws.Range("$col`3").Value2
    $row4Values[$col] = $ws.Range("$col`4").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col`3").Value2 = $row4Values[$col]
    $ws.Range("$col`4").Value2 = $row3Values[$col]
}
